$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-based values for rows 2-10, columns G-T
$updates = @{
    2  = @{ G=2.721294; H=8.163882000000001; I=0.03340571984979829; J=0.03340571984979829; M=0.8685706666666667; N=2.605712; O=0.2707495698024546; P=0.2707495698024546; Q=2.363636143776001; R=21.272725293984; S=0.009044584278274205; T=0.009044584278274205 }
    3  = @{ G=2.721294; H=8.163882000000001; I=0.03340571984979829; J=0.03340571984979829; O=0.4617398237148598; P=0.4617398237148598; Q=4.03097569887; R=36.27878128983; S=0.01542475119451386; T=0.01542475119451386 }
    4  = @{ G=2.721294; H=8.163882000000001; I=0.03340571984979829; J=0.03340571984979829; O=0.2675106064826855; P=0.2675106064826855; Q=2.33536008492; R=21.01824076428; S=0.008936384377010227; T=0.008936384377010227 }
    5  = @{ I=0.427684027063558; J=0.427684027063558; M=0.8685706666666667; N=2.605712; O=0.2707495698024546; P=0.2707495698024546; Q=30.26096815241067; R=272.348713371696; S=0.1157952663388397; T=0.1157952663388397 }
    6  = @{ I=0.427684027063558; J=0.427684027063558; O=0.4617398237148598; P=0.4617398237148598; S=0.1974787472619886; T=0.1974787472619886 }
    7  = @{ I=0.427684027063558; J=0.427684027063558; O=0.2675106064826855; P=0.2675106064826855; S=0.1144100134627297; T=0.1144100134627297 }
    8  = @{ I=0.5389102530866438; J=0.5389102530866438; M=0.8685706666666667; N=2.605712; O=0.2707495698024546; P=0.2707495698024546; Q=38.13082783949511; R=343.177450555456; S=0.1459097191853407; T=0.1459097191853407 }
    9  = @{ I=0.5389102530866438; J=0.5389102530866438; O=0.4617398237148598; P=0.4617398237148598; S=0.2488363252583574; T=0.2488363252583574 }
    10 = @{ I=0.5389102530866438; J=0.5389102530866438; O=0.2675106064826855; P=0.2675106064826855; S=0.1441642086429456; T=0.1441642086429456 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
